# Fix: add the missing "binomial" column ("algo"/"otro") to the `object`
# sheet and make that sheet the active tab/selection (previously
# `listed_simple` was active).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("object")

# New column F: header + values for rows 2-6
$ws.Range("F1").Value = "binomial"
$ws.Range("F2").Value = "algo"
$ws.Range("F3").Value = "otro"
$ws.Range("F4").Value = "algo"
$ws.Range("F5").Value = "algo"
$ws.Range("F6").Value = "otro"

# Make "object" the active sheet/tab and select F7 (just below the new data)
$ws.Activate() | Out-Null
$ws.Range("F7").Select() | Out-Null
